# ---------------------------------------------------------------------------
# Apply the commit's changes to buglist&newfeature.xlsx
#   1) MeetSDK: fix XOPlayer crash problem (support multi-instance player)
#      -> add two new bug rows (56 / 57) to the "bug" sheet
#   2) MeetPlayer: update crash update and breakpad util
#      -> update sheet view / selection state on "bug" and "newfeature"
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$bug = $wb.Worksheets.Item("bug")
$newfeature = $wb.Worksheets.Item("newfeature")

# ---------------------------------------------------------------------------
# 1) Add the two new bug rows for the XOPlayer crash / multi-instance fix
# ---------------------------------------------------------------------------

# Row 56 - high definition video stutter bug
$bug.Range("C56").Value = "N/A"
$bug.Range("D56").Value = "android"
$bug.Range("E56").Value = 20160420
$bug.Range("F56").Value = "XOPlayer 高清视频卡顿"
$bug.Range("G56").Value = " fixed"
$bug.Range("H56").Value = "read_sample proc 造成 video/audio 读取数据有耦合，无法较好同步。去除多余线程解决"

# Row 57 - XOPlayer does not support multiple instances
$bug.Range("C57").Value = "N/A"
$bug.Range("D57").Value = "android"
$bug.Range("E57").Value = 20160420
$bug.Range("F57").Value = "XOPlayer 不支持多实例"
$bug.Range("G57").Value = " fixed"
$bug.Range("H57").Value = "修改部分 static 变量"

# ---------------------------------------------------------------------------
# 2) Update sheet views / active selections
#    - "bug" becomes the active/selected tab with selection at G55
#    - "newfeature" loses the selected tab flag, selection moves to C45
#    (selecting a range on a sheet activates that sheet, so update
#     "newfeature" first and finish on "bug" so it ends up the active tab)
# ---------------------------------------------------------------------------

$newfeature.Activate()
$newfeature.Range("C45").Select()

$bug.Activate()
$bug.Range("G55").Select()
